# "Integrating Exel utility into the framework"
# Adds a small product/search/country test-data block in columns F:H
# of the existing worksheet, replacing the now-unused numeric sample
# data that used to live in row 4 (B4:E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (F1:H1)
$ws.Range("F1").Value = "product"
$ws.Range("G1").Value = "search"
$ws.Range("H1").Value = "country"

# Row 4 no longer holds the old numeric sample values (B4:E4) - it is
# repurposed to hold the new F:H sample data instead.
$ws.Range("B4:E4").ClearContents()

$ws.Range("F4").Value = "iphone X"
$ws.Range("G4").Value = "United"
$ws.Range("H4").Value = "United Kingdom"

# Size the new "country" column to fit its content.
$ws.Columns.Item(8).ColumnWidth = 14.6

# Match the author's final cursor position/selection.
$ws.Range("K11").Select()
